# [Kadastro App] Kayıtlar güncellendi - 16.07.2025 23:05:03
# The record that was in row 2 (Kayıt No 3 / İlçe / Tevhid / Gökhan ELGÜL)
# was removed; the remaining record (formerly row 3) shifts up to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

# Delete the entire second row (the "İlçe" / "Tevhid" / "Gökhan ELGÜL" record).
# This shifts the old row 3 up to become the new row 2, and the used range
# / dimension shrinks from A1:G3 to A1:G2 automatically.
$ws.Rows.Item(2).Delete()
